# Update NATMI LR-pair output sheet with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("M2").Value = 0.6068319999999999
$ws.Range("N2").Value = 1.820496
$ws.Range("O2").Value = 0.03392274820144286
$ws.Range("P2").Value = 0.03392274820144286
$ws.Range("Q2").Value = 5.556735946165333
$ws.Range("R2").Value = 50.010623515488
$ws.Range("S2").Value = 0.03288375181598083
$ws.Range("T2").Value = 0.03288375181598083

# Row 3
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("O3").Value = 0.4504903529585388
$ws.Range("P3").Value = 0.4504903529585388
$ws.Range("S3").Value = 0.4366925956061574
$ws.Range("T3").Value = 0.4366925956061574

# Row 4
$ws.Range("I4").Value = 0.9693716918425304
$ws.Range("J4").Value = 0.9693716918425304
$ws.Range("M4").Value = 9.223151
$ws.Range("N4").Value = 27.669453
$ws.Range("O4").Value = 0.5155868988400183
$ws.Range("P4").Value = 0.5155868988400183
$ws.Range("Q4").Value = 84.45601863219268
$ws.Range("R4").Value = 760.1041676897341
$ws.Range("S4").Value = 0.4997953444203921
$ws.Range("T4").Value = 0.4997953444203921

# Row 5
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("M5").Value = 0.6068319999999999
$ws.Range("N5").Value = 1.820496
$ws.Range("O5").Value = 0.03392274820144286
$ws.Range("P5").Value = 0.03392274820144286
$ws.Range("Q5").Value = 0.1755708592906667
$ws.Range("R5").Value = 1.580137733616
$ws.Range("S5").Value = 0.001038996385462041
$ws.Range("T5").Value = 0.00103899638546204

# Row 6
$ws.Range("G6").Value = 0.2893236666666667
$ws.Range("H6").Value = 0.867971
$ws.Range("I6").Value = 0.03062830815746963
$ws.Range("J6").Value = 0.03062830815746962
$ws.Range("O6").Value = 0.4504903529585388
$ws.Range("P6").Value = 0.4504903529585388
$ws.Range("Q6").Value = 2.331561638267333
$ws.Range("R6").Value = 20.984054744406
$ws.Range("S6").Value = 0.01379775735238139
$ws.Range("T6").Value = 0.01379775735238138

# Row 7
$ws.Range("G7").Value = 0.2893236666666667
$ws.Range("H7").Value = 0.867971
$ws.Range("I7").Value = 0.03062830815746963
$ws.Range("J7").Value = 0.03062830815746962
$ws.Range("M7").Value = 9.223151
$ws.Range("N7").Value = 27.669453
$ws.Range("O7").Value = 0.5155868988400183
$ws.Range("P7").Value = 0.5155868988400183
$ws.Range("Q7").Value = 2.668475865540334
$ws.Range("R7").Value = 24.016282789863
$ws.Range("S7").Value = 0.0157915544196262
$ws.Range("T7").Value = 0.0157915544196262
